$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename columns to snake_case machine-readable names
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Title-case the Spanish connector words (de/del/la/el/y/los) within
# municipality / state names throughout the data rows
$ws.Range("B4").Value = "Pabellón De Arteaga"
$ws.Range("B5").Value = "Rincón De Romos"
$ws.Range("B15").Value = "Amatenango De La Frontera"
$ws.Range("B36").Value = "Guadalupe Y Calvo"
$ws.Range("B38").Value = "Hidalgo Del Parral"
$ws.Range("A63").Value = "Ciudad De México"
$ws.Range("B88").Value = "San Juan Del Río"
$ws.Range("B89").Value = "San Luis Del Cordero"
$ws.Range("A94").Value = "Estado De México"
$ws.Range("B94").Value = "Almoloya De Alquisiras"
$ws.Range("B95").Value = "Almoloya De Juárez"
$ws.Range("B101").Value = "Ecatepec De Morelos"
$ws.Range("B105").Value = "Naucalpan De Juárez"
$ws.Range("B107").Value = "San Felipe Del Progreso"
$ws.Range("B111").Value = "Tlalnepantla De Baz"
$ws.Range("B123").Value = "Purísima Del Rincón"
$ws.Range("B126").Value = "Valle De Santiago"
$ws.Range("B129").Value = "Acapulco De Juárez"
$ws.Range("B130").Value = "Atoyac De Álvarez"
$ws.Range("B131").Value = "Chilpancingo De Los Bravo"
$ws.Range("B133").Value = "Coyuca De Benítez"
$ws.Range("B135").Value = "Zihuatanejo De Azueta"
$ws.Range("B141").Value = "Técpan De Galeana"
$ws.Range("B143").Value = "Tepecoacuilco De Trujano"
$ws.Range("B147").Value = "Atotonilco El Grande"
$ws.Range("B149").Value = "Mineral Del Chico"
$ws.Range("B150").Value = "Pachuca De Soto"
$ws.Range("B159").Value = "Autlán De Navarro"
$ws.Range("B162").Value = "Cuautitlán De García Barragán"
$ws.Range("B165").Value = "Encarnación De Díaz"
$ws.Range("B168").Value = "Lagos De Moreno"
$ws.Range("B171").Value = "San Juan De Los Lagos"
$ws.Range("B173").Value = "San Miguel El Alto"
$ws.Range("B175").Value = "Talpa De Allende"
$ws.Range("B177").Value = "Tepatitlán De Morelos"
$ws.Range("B178").Value = "Tizapán El Alto"
$ws.Range("B179").Value = "Tlajomulco De Zúñiga"
$ws.Range("B182").Value = "Unión De Tula"
$ws.Range("B186").Value = "Zapotlán El Grande"
$ws.Range("B191").Value = "Coalcomán De Vázquez Pallares"
$ws.Range("B223").Value = "Amatlán De Cañas"
$ws.Range("B228").Value = "Santa María Del Oro"
$ws.Range("B234").Value = "Ocotlán De Morelos"
$ws.Range("B235").Value = "Putla Villa De Guerrero"
$ws.Range("B249").Value = "Zimatlán De Álvarez"
$ws.Range("B260").Value = "Tepexi De Rodríguez"
$ws.Range("B262").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B269").Value = "Pinal De Amoles"
$ws.Range("B273").Value = "Santa María Del Río"
$ws.Range("B275").Value = "Villa De Ramos"
$ws.Range("B291").Value = "Nacozari De García"
$ws.Range("B310").Value = "Muñoz De Domingo Arenas"
$ws.Range("B312").Value = "Tepetitla De Lardizábal"
$ws.Range("B337").Value = "Nochistlán De Mejía"
$ws.Range("B342").Value = "Tlaltenango De Sánchez Román"

# Re-normalize a floating point rounding artifact in the percentage column
$ws.Range("D309").Value = 0.09988649262202044

# Remove the trailing footnote/metadata rows that are no longer part of
# the clean dataset. Delete bottom-most ranges first so row numbers for
# the other range stay valid.
$ws.Rows("476:480").Delete() | Out-Null
$ws.Rows("350:354").Delete() | Out-Null
